$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Designator lists: change the space-separated reference designators to
#    semicolon-separated ones (BOM export tool changed its delimiter).
#    Using Cells.Replace() edits the shared-string text without disturbing
#    any of the surrounding (untouched) shared strings.
# ---------------------------------------------------------------------------
$ws.Cells.Replace("B+1 B-1", "B+1;B-1")
$ws.Cells.Replace("P1 P3 P5 P6 P7 P8 P9 P10 P11 P12 P13 P14 P15 P16 P17 P18", "P1;P3;P5;P6;P7;P8;P9;P10;P11;P12;P13;P14;P15;P16;P17;P18")
$ws.Cells.Replace("P2 P4", "P2;P4")
$ws.Cells.Replace("U4 U5 U6 U7", "U4;U5;U6;U7")
$ws.Cells.Replace("C1 C7 C9 C12 C14 C16 C18", "C1;C7;C9;C12;C14;C16;C18")
$ws.Cells.Replace("C2 C8", "C2;C8")
$ws.Cells.Replace("C3 C5 C6 C10 C11 C13 C15 C17", "C3;C5;C6;C10;C11;C13;C15;C17")
$ws.Cells.Replace("R1 R2 R5 R6 R7 R8 R9 R10 R14 R15 R19 R20 R24 R25", "R1;R2;R5;R6;R7;R8;R9;R10;R14;R15;R19;R20;R24;R25")

# Row 17: R11 R16 R21 R26 @ 100k  ->  R11;R16;R21;R26 @ 10k
$ws.Cells.Replace("R11 R16 R21 R26", "R11;R16;R21;R26")
$ws.Cells.Replace("100k", "10k")

# Row 18 designator group shrinks from 8 refs to 4 (R12 R17 R22 R27), and its
# quantity drops from 8 to 4 to match. Value stays 2.4k.
$ws.Cells.Replace("R12 R13 R17 R18 R22 R23 R27 R28", "R12;R17;R22;R27")
$ws.Range("D18").Value = 4

# Row 19: the remaining refs (R13 R18 R23 R28) move here at 5.1k, replacing
# the old R29-R32 @ 1.2k row entirely.
$ws.Cells.Replace("R29 R30 R31 R32", "R13;R18;R23;R28")
$ws.Cells.Replace("1.2k", "5.1k")

# ---------------------------------------------------------------------------
# 2. Cell alignment: Id (A), Quantity (D) and Supplier-and-ref (E) columns
#    get explicit alignment for the data rows (2-19) as well as the header
#    row (1). Id/Supplier -> left, Quantity -> centered.
# ---------------------------------------------------------------------------
$ws.Range("A1:A19").HorizontalAlignment = -4131
$ws.Range("E1:E19").HorizontalAlignment = -4131
$ws.Range("D1:D19").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 3. Column widths were adjusted when the sheet was reformatted.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 53
$ws.Columns.Item(3).ColumnWidth = 32.833333333333336
$ws.Columns.Item(5).ColumnWidth = 20.166666666666668
$ws.Columns.Item(6).ColumnWidth = 19.666666666666668

# ---------------------------------------------------------------------------
# 4. Selection moved to C30 before the file was last saved.
# ---------------------------------------------------------------------------
$ws.Range("C30").Select()
